# Mise à jour de l'application
# Adds a new "training session" column (AT) dated 2025-09-12 (serial 45912)
# to the Présences sheet, with per-player attendance marks, mirroring the
# pattern of the existing AS column (dated 2025-09-10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell: new date in AT1, formatted/aligned like the other date
#     header cells (AS1 etc. use a center-aligned, vertically-centered
#     date number format). ---
$ws.Cells.Item(1, 46).Value = 45912
$ws.Cells.Item(1, 46).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(1, 46).HorizontalAlignment = -4108
$ws.Cells.Item(1, 46).VerticalAlignment = -4108

# --- Per-player attendance marks for the new session (column AT), one
#     letter code per row, matching the existing vocabulary used
#     throughout the sheet (P = Présent, R = Repos/Réserve, B = Blessure,
#     M = Malade). ---
$attendance = @(
    @{ Row = 2; Value = "P" },
    @{ Row = 3; Value = "R" },
    @{ Row = 4; Value = "P" },
    @{ Row = 5; Value = "B" },
    @{ Row = 6; Value = "B" },
    @{ Row = 7; Value = "P" },
    @{ Row = 8; Value = "P" },
    @{ Row = 9; Value = "P" },
    @{ Row = 10; Value = "B" },
    @{ Row = 11; Value = "P" },
    @{ Row = 12; Value = "P" },
    @{ Row = 13; Value = "B" },
    @{ Row = 14; Value = "P" },
    @{ Row = 15; Value = "P" },
    @{ Row = 16; Value = "P" },
    @{ Row = 17; Value = "P" },
    @{ Row = 18; Value = "P" },
    @{ Row = 19; Value = "P" },
    @{ Row = 20; Value = "P" },
    @{ Row = 21; Value = "R" },
    @{ Row = 22; Value = "P" },
    @{ Row = 23; Value = "R" },
    @{ Row = 24; Value = "P" },
    @{ Row = 25; Value = "P" },
    @{ Row = 26; Value = "P" },
    @{ Row = 27; Value = "M" },
    @{ Row = 28; Value = "P" },
    @{ Row = 29; Value = "P" }
)

foreach ($item in $attendance) {
    $cell = $ws.Cells.Item($item.Row, 46)
    $cell.Value = $item.Value
    $cell.HorizontalAlignment = -4108
}

# --- Recalculate the COUNTA/COUNTIF summary formulas (columns B-J) now
#     that the new column is inside their ranges. ---
$ws.Calculate()

# --- Restore the selection to mirror the sheet moving one column to the
#     right (AU27 -> AV27) now that a new last column exists. ---
$ws.Range("AV27").Select()
